# automate_finance.qmd needs to use the updated CVD files -> update cached
# CVD (cumulative-value-distribution) figures across the PES APAC location
# sheets so the "Commit/Forecast" monthly-weighting rows and a handful of
# "Manufacturing Voluntary Turnover" YTD figures reflect the refreshed data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: Bangkrang Nonthaburi
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Bangkrang Nonthaburi")

$ws.Range("L4").Value = 0.0128
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0

$ws.Range("E5").Value = 0.5
$ws.Range("E6").Value = 0.5

$ws.Range("E7").Value = 0.5
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 0.5
$ws.Range("N7").Value = 0.5
$ws.Range("O7").Value = 0.5
$ws.Range("P7").Value = 0.5
$ws.Range("Q7").Value = 0.5
$ws.Range("R7").Value = 0.5
$ws.Range("S7").Value = 0.5
$ws.Range("T7").Value = 0.5
$ws.Range("U7").Value = 0.5
$ws.Range("V7").Value = 0.5
$ws.Range("W7").Value = 0.5

$ws.Range("E8").Value = 0.0776
$ws.Range("E9").Value = 0.0776

$ws.Range("E10").Value = 0.0776
$ws.Range("L10").Value = 0.0148
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 0
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = 0
$ws.Range("Q10").Value = 0
$ws.Range("R10").Value = 0
$ws.Range("S10").Value = 0
$ws.Range("T10").Value = 0
$ws.Range("U10").Value = 0
$ws.Range("V10").Value = 0
$ws.Range("W10").Value = 0

# ---------------------------------------------------------------------
# Sheet: Yueyang China
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Yueyang China")

$ws.Range("L4").Value = 0.0133

$ws.Range("E7").Value = 0.0776
$ws.Range("E8").Value = 0.0776

$ws.Range("E9").Value = 0.0776
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 0
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = 0
$ws.Range("Q9").Value = 0
$ws.Range("R9").Value = 0
$ws.Range("S9").Value = 0
$ws.Range("T9").Value = 0
$ws.Range("U9").Value = 0
$ws.Range("V9").Value = 0
$ws.Range("W9").Value = 0

# ---------------------------------------------------------------------
# Sheet: Changzhou China Center (EPC) C
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Changzhou China Center (EPC) C")

$ws.Range("E7").Value = 0.0776
$ws.Range("E8").Value = 0.0776
$ws.Range("E9").Value = 0.0776

# ---------------------------------------------------------------------
# Sheet: Changzhou Epc China
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Changzhou Epc China")

$ws.Range("L7").ClearContents()

$ws.Range("E8").Value = 0.0776
$ws.Range("E9").Value = 0.0776

$ws.Range("E10").Value = 0.0776
$ws.Range("L10").Value = 0.0177
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 0
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = 0
$ws.Range("Q10").Value = 0
$ws.Range("R10").Value = 0
$ws.Range("S10").Value = 0
$ws.Range("T10").Value = 0
$ws.Range("U10").Value = 0
$ws.Range("V10").Value = 0
$ws.Range("W10").Value = 0

# ---------------------------------------------------------------------
# Sheet: Jiaxing China
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Jiaxing China")

$ws.Range("E7").Value = 0.0776
$ws.Range("E8").Value = 0.0776

$ws.Range("E9").Value = 0.0776
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 0
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = 0
$ws.Range("Q9").Value = 0
$ws.Range("R9").Value = 0
$ws.Range("S9").Value = 0
$ws.Range("T9").Value = 0
$ws.Range("U9").Value = 0
$ws.Range("V9").Value = 0
$ws.Range("W9").Value = 0

# ---------------------------------------------------------------------
# Sheet: Panyu Guangdong China
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Panyu Guangdong China")

$ws.Range("E5").Value = 0.0776
$ws.Range("E6").Value = 0.0776
$ws.Range("E7").Value = 0.0776

# ---------------------------------------------------------------------
# Sheet: Suzhou China
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Suzhou China")

$ws.Range("L4").Value = 0.0185

$ws.Range("E5").Value = 0.0776
$ws.Range("E6").Value = 0.0776

$ws.Range("E7").Value = 0.0776
$ws.Range("J7").Value = 0.0237
$ws.Range("K7").Value = 0.0169
$ws.Range("L7").Value = 0.0068
